# Updates the "Horarios Línea 141" workbook with the latest scrape
# (actualización 02:02:39) across all three sheets.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)   # LP1912
$ws2 = $wb.Worksheets.Item(2)   # LP1912-215
$ws3 = $wb.Worksheets.Item(3)   # 6203-6173

# ---------------------------------------------------------------
# Sheet 1: LP1912
# ---------------------------------------------------------------

$ws1.Cells.Item(2, 1).Value = "Última actualización: 02:02:39"
$ws1.Cells.Item(3, 1).Value = "Total filas: 3"

# Row 6 (was the 14_ABASTO / 6 min row, now the 15_ABASTO / 59 min row)
$ws1.Cells.Item(6, 1).Value = "02:02:39"
$ws1.Cells.Item(6, 2).Value = "03:01"
$ws1.Cells.Item(6, 3).Value = "15_ABASTO"
$ws1.Cells.Item(6, 4).Value = 59
$ws1.Cells.Item(6, 5).Value = "LP1912"

# Row 7 (was the 15_ABASTO / 102 min row, now the 215_ALUAR / 66 min row)
$ws1.Cells.Item(7, 1).Value = "02:02:39"
$ws1.Cells.Item(7, 2).Value = "03:08"
$ws1.Cells.Item(7, 3).Value = "215_ALUAR"
$ws1.Cells.Item(7, 4).Value = 66
$ws1.Cells.Item(7, 5).Value = "LP1912"

# Row 8 (new row)
$ws1.Cells.Item(8, 1).Value = "02:02:39"
$ws1.Cells.Item(8, 2).Value = "03:48"
$ws1.Cells.Item(8, 3).Value = "14_ABASTO"
$ws1.Cells.Item(8, 4).Value = 106
$ws1.Cells.Item(8, 5).Value = "LP1912"

# ---------------------------------------------------------------
# Sheet 2: LP1912-215
# ---------------------------------------------------------------

$ws2.Cells.Item(2, 1).Value = "Última actualización: 02:02:39"
$ws2.Cells.Item(3, 1).Value = "Total filas: 1"

# Bring the header row formatting over from sheet 1 (bold, bordered,
# centered) before filling in the header text / data row.
$ws1.Range("A5:E5").Copy()
$ws2.Range("A5:E5").PasteSpecial(-4122)  # xlPasteFormats

$ws2.Cells.Item(5, 1).Value = "Hora_Scrap"
$ws2.Cells.Item(5, 2).Value = "Hora_Llegada"
$ws2.Cells.Item(5, 3).Value = "Linea"
$ws2.Cells.Item(5, 4).Value = "Minutos"
$ws2.Cells.Item(5, 5).Value = "Parada"

$ws2.Cells.Item(6, 1).Value = "02:02:39"
$ws2.Cells.Item(6, 2).Value = "03:08"
$ws2.Cells.Item(6, 3).Value = "215_ALUAR"
$ws2.Cells.Item(6, 4).Value = 66
$ws2.Cells.Item(6, 5).Value = "LP1912"

# ---------------------------------------------------------------
# Sheet 3: 6203-6173
# ---------------------------------------------------------------

$ws3.Cells.Item(2, 1).Value = "Última actualización: 02:02:39"
